$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 6.2 Level of competition loses one answer option: "Niche circular market, no
# linear alternative available" (row 58, with its Example text in E58). Delete
# the whole row; everything below (6.3 Circular marketshare onward) shifts up.
$ws.Rows(58).Delete()

# Restore the user's on-screen selection/scroll position as left after the edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 51
$win.ScrollColumn = 1
[void]$ws.Range("C54").Select()
